# Add two new columns, I ("I0") and J ("IF"), to the right of the existing
# H ("IP") column. Header cells get the same style as the other headers
# (copied from H1 so we reuse the existing bold/bordered/centered style
# instead of minting a new one). I holds a constant 1 for every data row,
# J duplicates the H value for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy H1's formatting (bold font, thin border, centered alignment) onto
# the two new header cells so they match the existing header styling.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# Fill in the data rows: column I is always 1, column J mirrors column H.
for ($r = 2; $r -le 38; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
